$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "63.059.34"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.473.78"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("E4").Value = "  -0.37%  "
Set-TextValue "D5" "576.89"
$ws.Range("E5").Value = "  -0.36%  "
Set-TextValue "D6" "146.65"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "2.472.60"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  +1.30%  "
Set-TextValue "D12" "5.27"
$ws.Range("E12").Value = "  +0.45%  "
Set-TextValue "D13" "0.353"
$ws.Range("E13").Value = "  +0.38%  "
Set-TextValue "D14" "29.07"
$ws.Range("E14").Value = "  +7.87%  "
$ws.Range("D16").Value = "2.920.90"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "63.023.35"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "2.472.31"
$ws.Range("E18").Value = "  +1.13%  "
Set-TextValue "D19" "8.13"
Set-TextValue "D20" "11.04"
$ws.Range("E20").Value = "  +0.99%  "
Set-TextValue "D21" "329.48"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D22" "4.13"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D23" "2.22"
$ws.Range("E23").Value = "  +9.18%  "
$ws.Range("E24").Value = "  +0.04%  "
Set-TextValue "D25" "66.29"
$ws.Range("E25").Value = "  +0.82%  "
Set-TextValue "D26" "669.05"
$ws.Range("E26").Value = "  +8.73%  "
$ws.Range("E27").Value = "  +14.05%  "
$ws.Range("D28").Value = "0.0₃0989"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "2.592.91"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("E30").Value = "  -7.98%  "
Set-TextValue "D31" "1.45"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  +0.79%  "
Set-TextValue "D38" "153.24"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  -0.18%  "
Set-TextValue "D41" "18.76"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D45").Value = "0.0₆0308"
$ws.Range("E45").Value = "  +8.39%  "
Set-TextValue "D46" "150.74"
$ws.Range("E46").Value = "  +4.68%  "
Set-TextValue "D47" "15.15"
$ws.Range("E47").Value = "  +26.84%  "
$ws.Range("E48").Value = "  +0.77%  "
Set-TextValue "D49" "20.71"
$ws.Range("E49").Value = "  +2.79%  "
Set-TextValue "D50" "0.606"
$ws.Range("E50").Value = "  +1.42%  "
Set-TextValue "D51" "0.0514"
$ws.Range("E51").Value = "  -0.05%  "
